$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row for "TROPHIC SILICONE GEL" (product row 53) was removed from the
# report. All of the other product rows below it (54-67) move up one slot,
# but only in their data columns (C name, H balance, L order-limit, N price,
# P sell-price, Q transactions) - the running counter in column A keeps its
# original numbering. We reproduce that by copying each data column up by
# one row (pulling row 54's values into row 53, 55's into 54, etc.) using
# Copy + PasteSpecial (values and formats) so number formats/styles are
# preserved exactly.

$ws.Range("C54:C67").Copy()
$ws.Range("C53:C66").PasteSpecial(-4104)

$ws.Range("H54:H67").Copy()
$ws.Range("H53:H66").PasteSpecial(-4104)

$ws.Range("L54:L67").Copy()
$ws.Range("L53:L66").PasteSpecial(-4104)

$ws.Range("N54:N67").Copy()
$ws.Range("N53:N66").PasteSpecial(-4104)

$ws.Range("P54:P67").Copy()
$ws.Range("P53:P66").PasteSpecial(-4104)

$ws.Range("Q54:Q67").Copy()
$ws.Range("Q53:Q66").PasteSpecial(-4104)

# Row 67 no longer holds a product; it becomes the grand-total row (taking
# over what used to be row 68's role), so remove its old product formatting
# and merges, then write the new (reduced) total into P67:Q67.
$ws.Range("A67:B67").UnMerge()
$ws.Range("C67:G67").UnMerge()
$ws.Range("H67:K67").UnMerge()
$ws.Range("L67:M67").UnMerge()
$ws.Range("N67:O67").UnMerge()
$ws.Range("A67:O67").ClearContents()
$ws.Range("A67:O67").Style = "Normal"

$ws.Range("P68:Q68").Copy()
$ws.Range("P67:Q67").PasteSpecial(-4104)
$ws.Range("P67").Value2 = 4210.835
$ws.Range("P67:Q67").Merge()

# The footer (generated-by / timestamp / page) row used to be row 69; it is
# now row 68. Move its contents up one row and update the printed timestamp.
$ws.Range("A69:Q69").Copy()
$ws.Range("A68:Q68").PasteSpecial(-4104)
$ws.Range("G68:I68").Merge()
$ws.Range("K68:Q68").Merge()
$ws.Rows(68).RowHeight = $ws.Rows(69).RowHeight

$ws.Range("A68").Value2 = "Sunday, 1 June, 2025 5:55 PM"

# Remove the now-duplicated old footer row 69.
$ws.Rows(69).Delete()
